$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.889.58'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '1.709.16'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '312.76'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '0.3742'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("D8").Value = '49.39'
$ws.Range("E8").Value = '  +3.43%  '
$ws.Range("D9").Value = '0.3438'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").Value = '1.226'
$ws.Range("E10").Value = '  +4.70%  '
$ws.Range("D11").Value = '0.07551'
$ws.Range("E11").Value = '  +3.81%  '
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  +5.06%  '
$ws.Range("D14").Value = '6.324'
$ws.Range("E14").Value = '  +3.38%  '
$ws.Range("D15").Value = '7.095'
$ws.Range("D16").Value = '1.707.84'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("E17").Value = '  +2.25%  '
$ws.Range("D18").Value = '0.06737'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = '0.9986'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '84.07'
$ws.Range("E20").Value = '  +3.67%  '
$ws.Range("E21").Value = '  +5.35%  '
$ws.Range("D22").Value = '6.391'
$ws.Range("E22").Value = '  +4.46%  '
$ws.Range("D23").Value = '13.13'
$ws.Range("E23").Value = '  +7.64%  '
$ws.Range("D24").Value = '24.898.27'
$ws.Range("E24").Value = '  +2.14%  '
$ws.Range("D25").Value = '2.447'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").Value = '2.806'
$ws.Range("E26").Value = '  +5.40%  '
$ws.Range("D27").Value = '20.38'
$ws.Range("E27").Value = '  +4.37%  '
$ws.Range("D28").Value = '149.58'
$ws.Range("E28").Value = '  -2.61%  '
$ws.Range("D29").Value = '133.00'
$ws.Range("E29").Value = '  +4.25%  '
$ws.Range("D30").Value = '1.897.22'
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("D31").Value = '1.250'
$ws.Range("E31").Value = '  +28.25%  '
$ws.Range("D32").Value = '6.842'
$ws.Range("E32").Value = '  +8.26%  '
$ws.Range("D33").Value = '4.224'
$ws.Range("E33").Value = '  +4.07%  '
$ws.Range("D34").Value = '13.93'
$ws.Range("E34").Value = '  +12.55%  '
$ws.Range("D35").Value = '0.08803'
$ws.Range("E35").Value = '  +4.03%  '
$ws.Range("D36").Value = '1.770'
$ws.Range("E36").Value = '  +4.05%  '
$ws.Range("D37").Value = '5.627'
$ws.Range("E37").Value = '  +5.02%  '
$ws.Range("D38").Value = '0.06664'
$ws.Range("E38").Value = '  +2.61%  '
$ws.Range("D39").Value = '9.213'
$ws.Range("E39").Value = '  +3.76%  '
$ws.Range("E40").Value = '  +3.86%  '
$ws.Range("D41").Value = '0.2240'
$ws.Range("E41").Value = '  +6.60%  '
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").Value = '0.6472'
$ws.Range("E43").Value = '  +5.04%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").Value = '13.88'
$ws.Range("E45").Value = '  +5.03%  '
$ws.Range("D46").Value = '0.6164'
$ws.Range("E46").Value = '  +3.77%  '
$ws.Range("D47").Value = '3.839'
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("D48").Value = '2.133'
$ws.Range("E48").Value = '  +5.31%  '
$ws.Range("D49").Value = '130.02'
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("D50").Value = '0.07322'
$ws.Range("E50").Value = '  +1.81%  '
$ws.Range("D51").Value = '80.20'
$ws.Range("E51").Value = '  +5.63%  '
